# "temp manual sell order functionality"
# Applies a market-data refresh (new prices/pct_change/RSI, re-sorted by RSI)
# plus two new manual sell trades (PG, MRK) and a WMT buy, propagated through
# the stocks / portfolio / trades / summary sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) watchlist: refresh price / pct_change / rsi for every ticker and
#    re-sort the whole table by rsi ascending.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("watchlist")

$ws.Cells.Item(2,1).Value = "PFE"
$ws.Cells.Item(2,2).Value = 36.16
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = 34.67741935483868

$ws.Cells.Item(3,1).Value = "WMT"
$ws.Cells.Item(3,2).Value = 123.47
$ws.Cells.Item(3,3).Value = -0.3792
$ws.Cells.Item(3,4).Value = 36.87943262411344

$ws.Cells.Item(4,1).Value = "JNJ"
$ws.Cells.Item(4,2).Value = 148.65
$ws.Cells.Item(4,3).Value = 0.2698
$ws.Cells.Item(4,4).Value = 52.99455535390197

$ws.Cells.Item(5,1).Value = "INTC"
$ws.Cells.Item(5,2).Value = 61.93
$ws.Cells.Item(5,3).Value = -0.3059
$ws.Cells.Item(5,4).Value = 58.10276679841898

$ws.Cells.Item(6,1).Value = "MSFT"
$ws.Cells.Item(6,2).Value = 185.36
$ws.Cells.Item(6,3).Value = 0.2434
$ws.Cells.Item(6,4).Value = 59.15032679738569

$ws.Cells.Item(7,1).Value = "CSCO"
$ws.Cells.Item(7,2).Value = 46.94
$ws.Cells.Item(7,3).Value = 0.1066
$ws.Cells.Item(7,4).Value = 64.91677336747762

$ws.Cells.Item(8,1).Value = "DIS"
$ws.Cells.Item(8,2).Value = 122.18
$ws.Cells.Item(8,3).Value = 2.8884
$ws.Cells.Item(8,4).Value = 68.08707735062529

$ws.Cells.Item(9,1).Value = "VZ"
$ws.Cells.Item(9,2).Value = 56.83
$ws.Cells.Item(9,3).Value = 0.7624
$ws.Cells.Item(9,4).Value = 68.93353941267382

$ws.Cells.Item(10,1).Value = "CVX"
$ws.Cells.Item(10,2).Value = 97.18000000000001
$ws.Cells.Item(10,3).Value = 2.6296
$ws.Cells.Item(10,4).Value = 69.83210912906611

$ws.Cells.Item(11,1).Value = "V"
$ws.Cells.Item(11,2).Value = 196.87
$ws.Cells.Item(11,3).Value = 0.2597
$ws.Cells.Item(11,4).Value = 72.66143633071806

$ws.Cells.Item(12,1).Value = "MRK"
$ws.Cells.Item(12,2).Value = 82.06
$ws.Cells.Item(12,3).Value = 1.5092
$ws.Cells.Item(12,4).Value = 74.53754080522313

$ws.Cells.Item(13,1).Value = "IBM"
$ws.Cells.Item(13,2).Value = 129.05
$ws.Cells.Item(13,3).Value = 2.4206
$ws.Cells.Item(13,4).Value = 76.10250297973778

$ws.Cells.Item(14,1).Value = "WBA"
$ws.Cells.Item(14,2).Value = 43.61
$ws.Cells.Item(14,3).Value = 0.1378
$ws.Cells.Item(14,4).Value = 76.66666666666666

$ws.Cells.Item(15,1).Value = "RTX"
$ws.Cells.Item(15,2).Value = 67.47
$ws.Cells.Item(15,3).Value = 6.4196
$ws.Cells.Item(15,4).Value = 76.84507042253522

$ws.Cells.Item(16,1).Value = "AXP"
$ws.Cells.Item(16,2).Value = 105.4
$ws.Cells.Item(16,3).Value = 6.3894
$ws.Cells.Item(16,4).Value = 77.26999398677087

$ws.Cells.Item(17,1).Value = "AAPL"
$ws.Cells.Item(17,2).Value = 325.12
$ws.Cells.Item(17,3).Value = 0.5505
$ws.Cells.Item(17,4).Value = 77.62915129151303

$ws.Cells.Item(18,1).Value = "KO"
$ws.Cells.Item(18,2).Value = 47.9
$ws.Cells.Item(18,3).Value = 2.1322
$ws.Cells.Item(18,4).Value = 77.63157894736835

$ws.Cells.Item(19,1).Value = "GS"
$ws.Cells.Item(19,2).Value = 210.57
$ws.Cells.Item(19,3).Value = 3.1498
$ws.Cells.Item(19,4).Value = 77.93000990425882

$ws.Cells.Item(20,1).Value = "PG"
$ws.Cells.Item(20,2).Value = 118.53
$ws.Cells.Item(20,3).Value = 0.3981
$ws.Cells.Item(20,4).Value = 78.11634349030477

$ws.Cells.Item(21,1).Value = "JPM"
$ws.Cells.Item(21,2).Value = 104.27
$ws.Cells.Item(21,3).Value = 5.3978
$ws.Cells.Item(21,4).Value = 78.42565597667638

$ws.Cells.Item(22,1).Value = "XOM"
$ws.Cells.Item(22,2).Value = 49.24
$ws.Cells.Item(22,3).Value = 4.0795
$ws.Cells.Item(22,4).Value = 78.43347639484981

$ws.Cells.Item(23,1).Value = "CAT"
$ws.Cells.Item(23,2).Value = 127.09
$ws.Cells.Item(23,3).Value = 3.0571
$ws.Cells.Item(23,4).Value = 78.83211678832116

$ws.Cells.Item(24,1).Value = "MMM"
$ws.Cells.Item(24,2).Value = 161.21
$ws.Cells.Item(24,3).Value = 2.2647
$ws.Cells.Item(24,4).Value = 79.48823772183246

$ws.Cells.Item(25,1).Value = "HD"
$ws.Cells.Item(25,2).Value = 251
$ws.Cells.Item(25,3).ClearContents()
$ws.Cells.Item(25,4).Value = 83.98492699010842

$ws.Cells.Item(26,1).Value = "UNH"
$ws.Cells.Item(26,2).Value = 305.35
$ws.Cells.Item(26,3).Value = -0.3134
$ws.Cells.Item(26,4).Value = 86.22222222222234

$ws.Cells.Item(27,1).Value = "NKE"
$ws.Cells.Item(27,2).Value = 104.11
$ws.Cells.Item(27,3).Value = 3.3452
$ws.Cells.Item(27,4).Value = 88.48230353929206

$ws.Cells.Item(28,1).Value = "DOW"
$ws.Cells.Item(28,2).Value = 41.63
$ws.Cells.Item(28,3).Value = 2.5369
$ws.Cells.Item(28,4).Value = 89.59330143540673

$ws.Cells.Item(29,1).Value = "BA"
$ws.Cells.Item(29,2).Value = 173.16
$ws.Cells.Item(29,3).Value = 12.9476
$ws.Cells.Item(29,4).Value = 89.85074626865672

$ws.Cells.Item(30,1).Value = "MCD"
$ws.Cells.Item(30,2).Value = 193.29
$ws.Cells.Item(30,3).Value = 3.0385
$ws.Cells.Item(30,4).Value = 90.18867924528304

$ws.Cells.Item(31,1).Value = "TRV"
$ws.Cells.Item(31,2).Value = 114.47
$ws.Cells.Item(31,3).Value = 3.1261
$ws.Cells.Item(31,4).Value = 97.2792149866191

# ---------------------------------------------------------------------------
# 2) stocks: VZ / WMT refresh (WMT bought 5 more), JNJ refresh, PFE refresh;
#    PG and MRK are fully sold off -> their rows are removed.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("stocks")

$ws.Cells.Item(2,1).Value = "VZ"
$ws.Cells.Item(2,2).Value = 56.55
$ws.Cells.Item(2,3).Value = 56.83
$ws.Cells.Item(2,4).Value = 11
$ws.Cells.Item(2,5).Value = 625.13
$ws.Cells.Item(2,6).Value = 0.4951
$ws.Cells.Item(2,7).Value = 68.93353941267382
$ws.Cells.Item(2,8).Value = "buy"

$ws.Cells.Item(3,1).Value = "WMT"
$ws.Cells.Item(3,2).Value = 123.94
$ws.Cells.Item(3,3).Value = 123.47
$ws.Cells.Item(3,4).Value = 8
$ws.Cells.Item(3,5).Value = 987.76
$ws.Cells.Item(3,6).Value = -0.3792
$ws.Cells.Item(3,7).Value = 36.87943262411344
$ws.Cells.Item(3,8).Value = "buy"

$ws.Cells.Item(4,1).Value = "JNJ"
$ws.Cells.Item(4,2).Value = 149.11
$ws.Cells.Item(4,3).Value = 148.65
$ws.Cells.Item(4,4).Value = 3
$ws.Cells.Item(4,5).Value = 445.95
$ws.Cells.Item(4,6).Value = -0.3085
$ws.Cells.Item(4,7).Value = 52.99455535390197
$ws.Cells.Item(4,8).Value = "buy"

$ws.Cells.Item(5,1).Value = "PFE"
$ws.Cells.Item(5,2).Value = 35.46
$ws.Cells.Item(5,3).Value = 36.16
$ws.Cells.Item(5,4).Value = 23
$ws.Cells.Item(5,5).Value = 831.6799999999999
$ws.Cells.Item(5,6).Value = 1.9741
$ws.Cells.Item(5,7).Value = 34.67741935483868
$ws.Cells.Item(5,8).Value = "buy"

# PG (old row 4) and MRK (old row 5) have been sold out entirely - drop
# the now-stale rows 6:7 so the table shrinks back down to A1:H5.
$ws.Range("A6:H7").Delete()

# ---------------------------------------------------------------------------
# 3) portfolio: refreshed cash / stocks / total snapshot.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("portfolio")
$ws.Cells.Item(2,2).Value = 7733.720000000001
$ws.Cells.Item(3,2).Value = 2890.52
$ws.Cells.Item(4,2).Value = 10624.24

# ---------------------------------------------------------------------------
# 4) trades: append the new manual orders (WMT buy, PG sell, MRK sell).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("trades")

$ws.Cells.Item(39,1).Value = 37
$ws.Cells.Item(39,2).Value = "02/06/2020 18:50:13"
$ws.Cells.Item(39,3).Value = "WMT"
$ws.Cells.Item(39,4).Value = "buy"
$ws.Cells.Item(39,5).Value = 5
$ws.Cells.Item(39,6).Value = 619.7

$ws.Cells.Item(40,1).Value = 38
$ws.Cells.Item(40,2).Value = "03/06/2020 18:22:42"
$ws.Cells.Item(40,3).Value = "PG"
$ws.Cells.Item(40,4).Value = "sell"
$ws.Cells.Item(40,5).Value = 4
$ws.Cells.Item(40,6).Value = 474.12

$ws.Cells.Item(41,1).Value = 39
$ws.Cells.Item(41,2).Value = "03/06/2020 18:22:42"
$ws.Cells.Item(41,3).Value = "MRK"
$ws.Cells.Item(41,4).Value = "sell"
$ws.Cells.Item(41,5).Value = 6
$ws.Cells.Item(41,6).Value = 492.36

# column A carries the bold/border "row id" style throughout the sheet -
# replicate it on the three freshly appended rows.
$src = $ws.Range("A38")
$dst = $ws.Range("A39:A41")
$src.Copy()
$dst.PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5) summary: append the two new portfolio snapshots.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("summary")

$ws.Cells.Item(21,1).Value = "02/06/2020 18:50:11"
$ws.Cells.Item(21,2).Value = 6767.240000000002
$ws.Cells.Item(21,3).Value = 3845.63
$ws.Cells.Item(21,4).Value = 10612.87

$ws.Cells.Item(22,1).Value = "03/06/2020 18:22:41"
$ws.Cells.Item(22,2).Value = 7733.720000000001
$ws.Cells.Item(22,3).Value = 2890.52
$ws.Cells.Item(22,4).Value = 10624.24

$src = $ws.Range("A20")
$dst = $ws.Range("A21:A22")
$src.Copy()
$dst.PasteSpecial(-4122)

$excel.CutCopyMode = $false
